# Refresh computed market-board profit columns (H:N) on each crafting-job sheet.
# Source diff only touches raw <v> values (no formulas in this workbook), so
# we just re-stamp the updated numbers cell-by-cell per sheet/row.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 61: Mega-Potion of Strength
$ws.Range("H61").Value = 78005.08
$ws.Range("I61").Value = 804.9
$ws.Range("J61").Value = 335339
$ws.Range("K61").Value = 2414.7
$ws.Range("L61").Value = 1006017
$ws.Range("M61").Value = -2242.7
$ws.Range("N61").Value = -1006361

# Row 112: Superior Spiritbond Potion
$ws.Range("H112").Value = 19232848
$ws.Range("I112").Value = 3750
$ws.Range("J112").Value = 27779114
$ws.Range("K112").Value = 11250
$ws.Range("L112").Value = 83337342
$ws.Range("M112").Value = -10142
$ws.Range("N112").Value = -83339558

# Row 113: Starch Glue
$ws.Range("H113").Value = 4900.5
$ws.Range("I113").Value = 4980
$ws.Range("J113").Value = 4821
$ws.Range("K113").Value = 4980
$ws.Range("L113").Value = 4821
$ws.Range("M113").Value = -1726
$ws.Range("N113").Value = -11329

# Row 132: Growth Formula Lambda
$ws.Range("H132").Value = 3848278.8
$ws.Range("I132").Value = 4002000
$ws.Range("J132").Value = 5250
$ws.Range("K132").Value = 12006000
$ws.Range("L132").Value = 15750
$ws.Range("M132").Value = -12003470
$ws.Range("N132").Value = -20810

# Row 140: Book of Ra'Kaznar
$ws.Range("H140").Value = 30000
$ws.Range("J140").Value = 30000
$ws.Range("L140").Value = 30000
$ws.Range("N140").Value = -40360

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Steel Ingot
$ws.Range("H32").Value = 7605.43
$ws.Range("I32").Value = 7249.9194
$ws.Range("K32").Value = 7249.9194
$ws.Range("M32").Value = -6962.9194

# Row 74: Titanium Nugget
$ws.Range("H74").Value = 1710.683
$ws.Range("I74").Value = 1333.3715
$ws.Range("K74").Value = 1333.3715
$ws.Range("M74").Value = -459.3715

# Row 77: Titanium Nugget
$ws.Range("H77").Value = 1710.683
$ws.Range("I77").Value = 1333.3715
$ws.Range("K77").Value = 6666.8575
$ws.Range("M77").Value = -2298.8575

# Row 122: High Durium Nugget
$ws.Range("H122").Value = 2380.5938
$ws.Range("I122").Value = 1590.3636
$ws.Range("K122").Value = 4771.0908
$ws.Range("M122").Value = -2321.0908

# Row 132: Mountain Chromite Ingot
$ws.Range("H132").Value = 1966.3414
$ws.Range("I132").Value = 1406.1666
$ws.Range("J132").Value = 5999.6
$ws.Range("K132").Value = 4218.4998
$ws.Range("L132").Value = 17998.8
$ws.Range("M132").Value = -1688.4998
$ws.Range("N132").Value = -23058.8

# Row 134: Ruthenium Vambraces of Maiming
$ws.Range("H134").Value = 30702.857
$ws.Range("J134").Value = 30702.857
$ws.Range("L134").Value = 30702.857
$ws.Range("N134").Value = -40842.857

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Walnut Lumber
$ws.Range("H31").Value = 1719.65
$ws.Range("I31").Value = 1052.2373
$ws.Range("J31").Value = 2680.0732
$ws.Range("K31").Value = 1052.2373
$ws.Range("L31").Value = 2680.0732
$ws.Range("M31").Value = -757.2373
$ws.Range("N31").Value = -3270.0732

# Row 34: Walnut Lumber
$ws.Range("H34").Value = 1719.65
$ws.Range("I34").Value = 1052.2373
$ws.Range("J34").Value = 2680.0732
$ws.Range("K34").Value = 1052.2373
$ws.Range("L34").Value = 2680.0732
$ws.Range("M34").Value = -850.2373
$ws.Range("N34").Value = -3084.0732

# Row 62: Cedar Lumber
$ws.Range("H62").Value = 3784.8215
$ws.Range("I62").Value = 2386.1177
$ws.Range("J62").Value = 5946.4546
$ws.Range("K62").Value = 2386.1177
$ws.Range("L62").Value = 5946.4546
$ws.Range("M62").Value = -1762.1177
$ws.Range("N62").Value = -7194.4546

# Row 65: Cedar Lumber
$ws.Range("H65").Value = 3784.8215
$ws.Range("I65").Value = 2386.1177
$ws.Range("J65").Value = 5946.4546
$ws.Range("K65").Value = 11930.5885
$ws.Range("L65").Value = 29732.273
$ws.Range("M65").Value = -8810.588499999998
$ws.Range("N65").Value = -35972.273

# Row 132: Ginseng Lumber
$ws.Range("H132").Value = 2714.9333
$ws.Range("I132").Value = 2435
$ws.Range("J132").Value = 3484.75
$ws.Range("K132").Value = 7305
$ws.Range("L132").Value = 10454.25
$ws.Range("M132").Value = -4775
$ws.Range("N132").Value = -15514.25

$ws = $wb.Worksheets.Item("CUL")
# Row 82: Baked Pipira Pira
$ws.Range("H82").Value = 3500

# Row 85: Baked Pipira Pira
$ws.Range("H85").Value = 3500

# Row 97: Cottonseed Oil
$ws.Range("H97").Value = 1460
$ws.Range("I97").Value = 433.33334
$ws.Range("K97").Value = 1300.00002
$ws.Range("M97").Value = -804.0000199999999

# Row 126: Glory Be Soup
$ws.Range("H126").Value = 1155
$ws.Range("I126").Value = 482.5
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 1447.5
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = 3492.5
$ws.Range("N126").Value = -17380

# Row 132: Cooking Mezcal
$ws.Range("H132").Value = 1955.8214
$ws.Range("J132").Value = 2281.353
$ws.Range("L132").Value = 20532.177
$ws.Range("N132").Value = -25592.177

$ws = $wb.Worksheets.Item("GSM")
# Row 122: Ametrine
$ws.Range("H122").Value = 4187.0376
$ws.Range("I122").Value = 3328.6875
$ws.Range("J122").Value = 5495
$ws.Range("K122").Value = 9986.0625
$ws.Range("L122").Value = 16485
$ws.Range("M122").Value = -7536.0625
$ws.Range("N122").Value = -21385

# Row 126: Phrygian Gold Ingot
$ws.Range("H126").Value = 628676.9
$ws.Range("I126").Value = 2402.2856
$ws.Range("J126").Value = 1115779.4
$ws.Range("K126").Value = 7206.8568
$ws.Range("L126").Value = 3347338.2
$ws.Range("M126").Value = -4736.8568
$ws.Range("N126").Value = -3352278.2

# Row 132: Lar Ingot
$ws.Range("H132").Value = 3373.8293
$ws.Range("I132").Value = 3615.8635
$ws.Range("K132").Value = 10847.5905
$ws.Range("M132").Value = -8317.5905

$ws = $wb.Worksheets.Item("LTW")
# Row 17: Hard Leather Harness
$ws.Range("H17").Value = 21007.363
$ws.Range("J17").Value = 23008.1
$ws.Range("L17").Value = 23008.1
$ws.Range("N17").Value = -23348.1

# Row 94: Gaganaskin Hat of Aiming
$ws.Range("H94").Value = 30776.666
$ws.Range("J94").Value = 30776.666
$ws.Range("L94").Value = 30776.666
$ws.Range("N94").Value = -32128.666

# Row 106: Gazelleskin Boots of Casting
$ws.Range("H106").Value = 27857.143
$ws.Range("J106").Value = 27857.143
$ws.Range("L106").Value = 27857.143
$ws.Range("N106").Value = -30381.143

# Row 132: Silver Lobo Leather
$ws.Range("H132").Value = 2311.54
$ws.Range("I132").Value = 1644.6177
$ws.Range("J132").Value = 3728.75
$ws.Range("K132").Value = 4933.8531
$ws.Range("L132").Value = 11186.25
$ws.Range("M132").Value = -2403.8531
$ws.Range("N132").Value = -16246.25

$ws = $wb.Worksheets.Item("WVR")
# Row 13: Hempen Acton
$ws.Range("H13").Value = 72670.664
$ws.Range("I13").Value = 18000
$ws.Range("K13").Value = 18000
$ws.Range("M13").Value = -17860

# Row 41: Linen Halfgloves
$ws.Range("H41").Value = 8610.833000000001
$ws.Range("J41").Value = 8610.833000000001
$ws.Range("L41").Value = 8610.833000000001
$ws.Range("N41").Value = -9390.833000000001

# Row 101: Serge Hose of Aiming
$ws.Range("H101").Value = 20875
$ws.Range("J101").Value = 20875
$ws.Range("L101").Value = 20875
$ws.Range("N101").Value = -27365

# Row 107: Bright Linen Yarn
$ws.Range("H107").Value = 990.2143
$ws.Range("I107").Value = 306.3
$ws.Range("J107").Value = 2700
$ws.Range("K107").Value = 918.9000000000001
$ws.Range("L107").Value = 8100
$ws.Range("M107").Value = 1001.1
$ws.Range("N107").Value = -11940

# Row 132: Snow Cotton Cloth
$ws.Range("H132").Value = 12516.125
$ws.Range("I132").Value = 2048.2942
$ws.Range("J132").Value = 71833.836
$ws.Range("K132").Value = 6144.882599999999
$ws.Range("L132").Value = 215501.508
$ws.Range("M132").Value = -3614.882599999999
$ws.Range("N132").Value = -220561.508

# Row 138: Rroneek Serge Halfgloves of Healing
$ws.Range("H138").Value = 29700
$ws.Range("J138").Value = 29700
$ws.Range("L138").Value = 29700
$ws.Range("N138").Value = -39980
